# urbanicity back to PCA formulation
#
# The sheet previously listed four separate "urbanicity" covariates
# (built_population_2014, nightlights_composite, all_population_count_2015,
# travel_times_2015) in rows 19-21 (plus B18). This reverts to a single
# PCA-derived "urbanscore" covariate:
#   - delete rows 19:21 (the three extra covariate rows)
#   - rename the remaining covariate cell (B18) to urbanscore_cont_scale_clst
#   - tidy up the selection to the single remaining cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three now-obsolete covariate rows; everything below shifts up
# (old row 22 "I:Pv18s" / "Pv18s" becomes the new row 19).
$ws.Range("A19:B21").EntireRow.Delete() | Out-Null

# Collapse the four urbanicity covariates into the single PCA-based one.
$ws.Range("B18").Value = "urbanscore_cont_scale_clst"

# Match the saved selection state (single cell instead of the old B18:B21).
$ws.Range("B18").Select() | Out-Null
